$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 2.291541108688389
$ws.Range("D2").Value = 4.050898842812317
$ws.Range("E2").Value = 16.70908072926618
$ws.Range("F2").Value = 31.60504126112289
$ws.Range("G2").Value = 46.55409832425982
$ws.Range("H2").Value = 15.40832056857072
$ws.Range("I2").Value = 24.58747617932364
$ws.Range("N2").Value = 18.99769123883711
$ws.Range("C3").Value = 2.242653789023274
$ws.Range("D3").Value = 3.966652200631139
$ws.Range("E3").Value = 15.73539356624581
$ws.Range("F3").Value = 30.3134804677301
$ws.Range("G3").Value = 44.04312956734825
$ws.Range("H3").Value = 15.06860057062744
$ws.Range("I3").Value = 23.73397333187833
$ws.Range("N3").Value = 18.40031508502701
$ws.Range("C4").Value = 2.211807746959753
$ws.Range("D4").Value = 3.915703165446748
$ws.Range("E4").Value = 15.11294793995385
$ws.Range("F4").Value = 29.50471783009843
$ws.Range("G4").Value = 42.44369329197139
$ws.Range("H4").Value = 14.86204866385766
$ws.Range("I4").Value = 23.20411622102332
$ws.Range("N4").Value = 18.02485520896362
$ws.Range("C5").Value = 2.199035036526418
$ws.Range("D5").Value = 3.895165308247801
$ws.Range("E5").Value = 14.85337459359302
$ws.Range("F5").Value = 29.17174623420898
$ws.Range("G5").Value = 41.77824074587323
$ws.Range("H5").Value = 14.77851342090213
$ws.Range("I5").Value = 22.98711071101861
$ws.Range("N5").Value = 17.86990355188116
$ws.Range("C6").Value = 2.196902070549409
$ws.Range("D6").Value = 3.891769386998015
$ws.Range("E6").Value = 14.80992375570737
$ws.Range("F6").Value = 29.1162691702013
$ws.Range("G6").Value = 41.66694468322564
$ws.Range("H6").Value = 14.76468452929513
$ws.Range("I6").Value = 22.95102309836007
$ws.Range("N6").Value = 17.84406337566597
$ws.Range("C7").Value = 2.211636302568591
$ws.Range("D7").Value = 3.915425239108351
$ws.Range("E7").Value = 15.10947083176046
$ws.Range("F7").Value = 29.50024021401694
$ws.Range("G7").Value = 42.4347728965714
$ws.Range("H7").Value = 14.86091933738302
$ws.Range("I7").Value = 23.20119348434856
$ws.Range("N7").Value = 18.02277304766463
$ws.Range("C8").Value = 2.274860173089194
$ws.Range("D8").Value = 4.021706081111173
$ws.Range("E8").Value = 16.37860076599894
$ws.Range("F8").Value = 31.16329760214364
$ws.Range("G8").Value = 45.70074826798113
$ws.Range("H8").Value = 15.29083622336184
$ws.Range("I8").Value = 24.29460900471237
$ws.Range("N8").Value = 18.79364780656866
$ws.Range("C9").Value = 2.392031015747773
$ws.Range("D9").Value = 4.235053467491193
$ws.Range("E9").Value = 18.82411019586357
$ws.Range("F9").Value = 34.27788795284263
$ws.Range("G9").Value = 51.61769986820902
$ws.Range("H9").Value = 16.14471258504169
$ws.Range("I9").Value = 26.37774005021591
$ws.Range("N9").Value = 20.2273683202997
$ws.Range("C10").Value = 2.473682401355542
$ws.Range("D10").Value = 4.393165847162185
$ws.Range("E10").Value = 20.53581442695382
$ws.Range("F10").Value = 36.4516128695236
$ws.Range("G10").Value = 55.6358714946854
$ws.Range("H10").Value = 16.77170479421687
$ws.Range("I10").Value = 27.85301576224247
$ws.Range("N10").Value = 21.22223697909767
$ws.Range("C11").Value = 2.509813349873522
$ws.Range("D11").Value = 4.465060646647207
$ws.Range("E11").Value = 21.2729129993351
$ws.Range("F11").Value = 37.41136784536977
$ws.Range("G11").Value = 57.38764912677016
$ws.Range("H11").Value = 17.05548680334735
$ws.Range("I11").Value = 28.50897267896941
$ws.Range("N11").Value = 21.66018057919901
$ws.Range("C12").Value = 2.523345721851969
$ws.Range("D12").Value = 4.492254776971405
$ws.Range("E12").Value = 21.54610482686304
$ws.Range("F12").Value = 37.77032924410633
$ws.Range("G12").Value = 58.03976241435712
$ws.Range("H12").Value = 17.16263380323386
$ws.Range("I12").Value = 28.75496003551665
$ws.Range("N12").Value = 21.82377585682186
$ws.Range("C13").Value = 2.520437998086978
$ws.Range("D13").Value = 4.486399811963222
$ws.Range("E13").Value = 21.48753115452632
$ws.Range("F13").Value = 37.69322384252373
$ws.Range("G13").Value = 57.89982250672377
$ws.Range("H13").Value = 17.13957344167379
$ws.Range("I13").Value = 28.70209278663542
$ws.Range("N13").Value = 21.78864458690801
$ws.Range("C14").Value = 2.510929686690126
$ws.Range("D14").Value = 4.467298688734616
$ws.Range("E14").Value = 21.2955072684879
$ws.Range("F14").Value = 37.44099107847981
$ws.Range("G14").Value = 57.44152579165135
$ws.Range("H14").Value = 17.0643087612632
$ws.Range("I14").Value = 28.52925971419292
$ws.Range("N14").Value = 21.67368539489659
$ws.Range("C15").Value = 2.505085972457943
$ws.Range("D15").Value = 4.455593920829401
$ws.Range("E15").Value = 21.17711600208215
$ws.Range("F15").Value = 37.28590017163895
$ws.Range("G15").Value = 57.15933325238638
$ws.Range("H15").Value = 17.01816284569774
$ws.Range("I15").Value = 28.42307434073396
$ws.Range("N15").Value = 21.60297336126124
$ws.Range("C16").Value = 2.471300429075009
$ws.Range("D16").Value = 4.388464436848404
$ws.Range("E16").Value = 20.48681080823547
$ws.Range("F16").Value = 36.38827920269986
$ws.Range("G16").Value = 55.51982985930501
$ws.Range("H16").Value = 16.75312061138612
$ws.Range("I16").Value = 27.80982219642646
$ws.Range("N16").Value = 21.19330956972086
$ws.Range("C17").Value = 2.450311739149582
$ws.Range("D17").Value = 4.347255452497016
$ws.Range("E17").Value = 20.05271765595111
$ws.Range("F17").Value = 35.82994522078612
$ws.Range("G17").Value = 54.49431812167074
$ws.Range("H17").Value = 16.59007803468783
$ws.Range("I17").Value = 27.42955806999067
$ws.Range("N17").Value = 20.93814219015166
$ws.Range("C18").Value = 2.438144496188168
$ws.Range("D18").Value = 4.323551681054679
$ws.Range("E18").Value = 19.79912380347964
$ws.Range("F18").Value = 35.5060835398215
$ws.Range("G18").Value = 53.89732073470162
$ws.Range("H18").Value = 16.49617078439743
$ws.Range("I18").Value = 27.20942584592137
$ws.Range("N18").Value = 20.79000725568362
$ws.Range("C19").Value = 2.434008671566446
$ws.Range("D19").Value = 4.315526567921878
$ws.Range("E19").Value = 19.71258750246551
$ws.Range("F19").Value = 35.39597160121749
$ws.Range("G19").Value = 53.69396992151733
$ws.Range("H19").Value = 16.46435656607968
$ws.Range("I19").Value = 27.13465741008431
$ws.Range("N19").Value = 20.73962067985785
$ws.Range("C20").Value = 2.452555894226304
$ws.Range("D20").Value = 4.351642565151272
$ws.Range("E20").Value = 20.09933246648837
$ws.Range("F20").Value = 35.88966489208477
$ws.Range("G20").Value = 54.60422800683271
$ws.Range("H20").Value = 16.60744839901975
$ws.Range("I20").Value = 27.47018597133863
$ws.Range("N20").Value = 20.96544799483449
$ws.Range("C21").Value = 2.513726599247696
$ws.Range("D21").Value = 4.472910192488118
$ws.Range("E21").Value = 21.35206995216972
$ws.Range("F21").Value = 37.51520159805676
$ws.Range("G21").Value = 57.57644598078067
$ws.Range("H21").Value = 17.08642520133454
$ws.Range("I21").Value = 28.58009206386674
$ws.Range("N21").Value = 21.70751365554064
$ws.Range("C22").Value = 2.552831431934818
$ws.Range("D22").Value = 4.557067552216613
$ws.Range("E22").Value = 22.13625475371152
$ws.Range("F22").Value = 38.55138096152353
$ws.Range("G22").Value = 59.45330909550263
$ws.Range("H22").Value = 17.39758237128472
$ws.Range("I22").Value = 29.29134882595776
$ws.Range("N22").Value = 22.17935961385674
$ws.Range("C23").Value = 2.532041617467994
$ws.Range("D23").Value = 4.509802200768107
$ws.Range("E23").Value = 21.72086697686762
$ws.Range("F23").Value = 38.00083668319122
$ws.Range("G23").Value = 58.45768323017286
$ws.Range("H23").Value = 17.23171805887528
$ws.Range("I23").Value = 28.91309976745497
$ws.Range("N23").Value = 21.92877110911181
$ws.Range("C24").Value = 2.451541624646266
$ws.Range("D24").Value = 4.349659187433638
$ws.Range("E24").Value = 20.07827046014655
$ws.Range("F24").Value = 35.8626745391651
$ws.Range("G24").Value = 54.5545608265352
$ws.Range("H24").Value = 16.59959578446747
$ws.Range("I24").Value = 27.45182279300356
$ws.Range("N24").Value = 20.95310750188673
$ws.Range("C25").Value = 2.361092962591343
$ws.Range("D25").Value = 4.176981437835304
$ws.Range("E25").Value = 18.15660052169967
$ws.Range("F25").Value = 33.45381830685568
$ws.Range("G25").Value = 50.07289128336305
$ws.Range("H25").Value = 15.91329210517361
$ws.Range("I25").Value = 25.82270765140368
$ws.Range("N25").Value = 19.84905939529497
